# Regenerate save_data "K" column (strikeouts, formerly "Strike#") with
# freshly calculated s_vals for each outing row (rows 2-55).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 0
    4  = 1
    5  = 3
    6  = 0
    7  = 2
    8  = 2
    9  = 1
    10 = 0
    11 = 1
    12 = 1
    13 = 3
    14 = 1
    15 = 0
    16 = 2
    17 = 1
    18 = 1
    19 = 2
    20 = 0
    21 = 2
    22 = 1
    23 = 2
    24 = 1
    25 = 0
    26 = 0
    27 = 3
    28 = 0
    29 = 3
    30 = 1
    31 = 0
    32 = 0
    33 = 1
    34 = 1
    35 = 0
    36 = 2
    37 = 2
    38 = 1
    39 = 0
    40 = 2
    41 = 2
    42 = 0
    43 = 0
    44 = 1
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = 2
    50 = 0
    51 = 0
    52 = 1
    53 = 0
    54 = 1
    55 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
